$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.46510533333333
$ws.Range("H2").Value = 31.395316
$ws.Range("I2").Value = 0.5554075997074865
$ws.Range("J2").Value = 0.5554075997074865
$ws.Range("M2").Value = 0.103005
$ws.Range("N2").Value = 0.309015
$ws.Range("O2").Value = 0.004680825815734043
$ws.Range("P2").Value = 0.004680825815734042
$ws.Range("Q2").Value = 1.07795817486
$ws.Range("R2").Value = 9.701623573739999
$ws.Range("S2").Value = 0.002599766230965683
$ws.Range("T2").Value = 0.002599766230965682
$ws.Range("G3").Value = 10.46510533333333
$ws.Range("H3").Value = 31.395316
$ws.Range("I3").Value = 0.5554075997074865
$ws.Range("J3").Value = 0.5554075997074865
$ws.Range("N3").Value = 54.36902
$ws.Range("O3").Value = 0.8235584434158876
$ws.Range("P3").Value = 0.8235584434158876
$ws.Range("Q3").Value = 189.6591737233689
$ws.Range("R3").Value = 1706.93256351032
$ws.Range("S3").Value = 0.457410618276452
$ws.Range("T3").Value = 0.457410618276452
$ws.Range("G4").Value = 10.46510533333333
$ws.Range("H4").Value = 31.395316
$ws.Range("I4").Value = 0.5554075997074865
$ws.Range("J4").Value = 0.5554075997074865
$ws.Range("M4").Value = 3.761887333333334
$ws.Range("N4").Value = 11.285662
$ws.Range("O4").Value = 0.1709503358647596
$ws.Range("P4").Value = 0.1709503358647596
$ws.Range("Q4").Value = 39.36854719546578
$ws.Range("R4").Value = 354.316924759192
$ws.Range("S4").Value = 0.0949471157118348
$ws.Range("T4").Value = 0.09494711571183478
$ws.Range("G5").Value = 10.46510533333333
$ws.Range("H5").Value = 31.395316
$ws.Range("I5").Value = 0.5554075997074865
$ws.Range("J5").Value = 0.5554075997074865
$ws.Range("M5").Value = 0.01783333333333333
$ws.Range("N5").Value = 0.0535
$ws.Range("O5").Value = 0.0008103949036188254
$ws.Range("P5").Value = 0.0008103949036188253
$ws.Range("Q5").Value = 0.1866277117777778
$ws.Range("R5").Value = 1.679649406
$ws.Range("S5").Value = 0.0004500994882341117
$ws.Range("T5").Value = 0.0004500994882341116
$ws.Range("I6").Value = 0.3053945925621632
$ws.Range("J6").Value = 0.3053945925621632
$ws.Range("M6").Value = 0.103005
$ws.Range("N6").Value = 0.309015
$ws.Range("O6").Value = 0.004680825815734043
$ws.Range("P6").Value = 0.004680825815734042
$ws.Range("Q6").Value = 0.5927225298750001
$ws.Range("R6").Value = 5.334502768875001
$ws.Range("S6").Value = 0.001429498892850553
$ws.Range("T6").Value = 0.001429498892850553
$ws.Range("I7").Value = 0.3053945925621632
$ws.Range("J7").Value = 0.3053945925621632
$ws.Range("N7").Value = 54.36902
$ws.Range("O7").Value = 0.8235584434158876
$ws.Range("P7").Value = 0.8235584434158876
$ws.Range("R7").Value = 938.5683145835002
$ws.Range("S7").Value = 0.2515102952781243
$ws.Range("T7").Value = 0.2515102952781243
$ws.Range("I8").Value = 0.3053945925621632
$ws.Range("J8").Value = 0.3053945925621632
$ws.Range("M8").Value = 3.761887333333334
$ws.Range("N8").Value = 11.285662
$ws.Range("O8").Value = 0.1709503358647596
$ws.Range("P8").Value = 0.1709503358647596
$ws.Range("Q8").Value = 21.64705963126111
$ws.Range("R8").Value = 194.82353668135
$ws.Range("S8").Value = 0.05220730816978322
$ws.Range("T8").Value = 0.05220730816978321
$ws.Range("I9").Value = 0.3053945925621632
$ws.Range("J9").Value = 0.3053945925621632
$ws.Range("M9").Value = 0.01783333333333333
$ws.Range("N9").Value = 0.0535
$ws.Range("O9").Value = 0.0008103949036188254
$ws.Range("P9").Value = 0.0008103949036188253
$ws.Range("Q9").Value = 0.1026184986111111
$ws.Range("R9").Value = 0.9235664875000001
$ws.Range("S9").Value = 0.0002474902214051247
$ws.Range("T9").Value = 0.0002474902214051247
$ws.Range("G10").Value = 2.146766
$ws.Range("H10").Value = 6.440298
$ws.Range("I10").Value = 0.1139338891693565
$ws.Range("J10").Value = 0.1139338891693565
$ws.Range("M10").Value = 0.103005
$ws.Range("N10").Value = 0.309015
$ws.Range("O10").Value = 0.004680825815734043
$ws.Range("P10").Value = 0.004680825815734042
$ws.Range("Q10").Value = 0.22112763183
$ws.Range("R10").Value = 1.99014868647
$ws.Range("S10").Value = 0.0005333046897109053
$ws.Range("T10").Value = 0.0005333046897109053
$ws.Range("G11").Value = 2.146766
$ws.Range("H11").Value = 6.440298
$ws.Range("I11").Value = 0.1139338891693565
$ws.Range("J11").Value = 0.1139338891693565
$ws.Range("N11").Value = 54.36902
$ws.Range("O11").Value = 0.8235584434158876
$ws.Range("P11").Value = 0.8235584434158876
$ws.Range("Q11").Value = 38.90585452977333
$ws.Range("R11").Value = 350.15269076796
$ws.Range("S11").Value = 0.09383121641663351
$ws.Range("T11").Value = 0.09383121641663353
$ws.Range("G12").Value = 2.146766
$ws.Range("H12").Value = 6.440298
$ws.Range("I12").Value = 0.1139338891693565
$ws.Range("J12").Value = 0.1139338891693565
$ws.Range("M12").Value = 3.761887333333334
$ws.Range("N12").Value = 11.285662
$ws.Range("O12").Value = 0.1709503358647596
$ws.Range("P12").Value = 0.1709503358647596
$ws.Range("Q12").Value = 8.075891823030666
$ws.Range("R12").Value = 72.683026407276
$ws.Range("S12").Value = 0.0194770366198798
$ws.Range("T12").Value = 0.0194770366198798
$ws.Range("G13").Value = 2.146766
$ws.Range("H13").Value = 6.440298
$ws.Range("I13").Value = 0.1139338891693565
$ws.Range("J13").Value = 0.1139338891693565
$ws.Range("M13").Value = 0.01783333333333333
$ws.Range("N13").Value = 0.0535
$ws.Range("O13").Value = 0.0008103949036188254
$ws.Range("P13").Value = 0.0008103949036188253
$ws.Range("Q13").Value = 0.03828399366666666
$ws.Range("R13").Value = 0.344555943
$ws.Range("S13").Value = 0.00009233144313231861
$ws.Range("T13").Value = 0.00009233144313231861
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4760280000000001
$ws.Range("H14").Value = 1.428084
$ws.Range("I14").Value = 0.02526391856099382
$ws.Range("J14").Value = 0.02526391856099382
$ws.Range("M14").Value = 0.103005
$ws.Range("N14").Value = 0.309015
$ws.Range("O14").Value = 0.004680825815734043
$ws.Range("P14").Value = 0.004680825815734042
$ws.Range("Q14").Value = 0.04903326414
$ws.Range("R14").Value = 0.44129937726
$ws.Range("S14").Value = 0.0001182560022069023
$ws.Range("T14").Value = 0.0001182560022069023
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4760280000000001
$ws.Range("H15").Value = 1.428084
$ws.Range("I15").Value = 0.02526391856099382
$ws.Range("J15").Value = 0.02526391856099382
$ws.Range("N15").Value = 54.36902
$ws.Range("O15").Value = 0.8235584434158876
$ws.Range("P15").Value = 0.8235584434158876
$ws.Range("Q15").Value = 8.627058617520001
$ws.Range("R15").Value = 77.64352755768
$ws.Range("S15").Value = 0.02080631344467782
$ws.Range("T15").Value = 0.02080631344467782
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4760280000000001
$ws.Range("H16").Value = 1.428084
$ws.Range("I16").Value = 0.02526391856099382
$ws.Range("J16").Value = 0.02526391856099382
$ws.Range("M16").Value = 3.761887333333334
$ws.Range("N16").Value = 11.285662
$ws.Range("O16").Value = 0.1709503358647596
$ws.Range("P16").Value = 0.1709503358647596
$ws.Range("Q16").Value = 1.790763703512
$ws.Range("R16").Value = 16.116873331608
$ws.Range("S16").Value = 0.004318875363261828
$ws.Range("T16").Value = 0.004318875363261828
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4760280000000001
$ws.Range("H17").Value = 1.428084
$ws.Range("I17").Value = 0.02526391856099382
$ws.Range("J17").Value = 0.02526391856099382
$ws.Range("M17").Value = 0.01783333333333333
$ws.Range("N17").Value = 0.0535
$ws.Range("O17").Value = 0.0008103949036188254
$ws.Range("P17").Value = 0.0008103949036188253
$ws.Range("Q17").Value = 0.008489166000000001
$ws.Range("R17").Value = 0.076402494
$ws.Range("S17").Value = 0.00002047375084727044
$ws.Range("T17").Value = 0.00002047375084727044
